$d = $word.ActiveDocument

function Find-ParaByText($needle) {
    $p = $d.Paragraphs.First
    while ($p -ne $null) {
        $t = $p.Range.Text
        if ($t -ne $null -and $t.Contains($needle)) {
            return $p
        }
        $p = $p.Next()
    }
    return $null
}

# ---------------------------------------------------------------------
# Row "Week beginning April 5th": 3 trailing empty paragraphs become 2
# paragraphs of notes text (delete the 3rd, fill the first two).
# ---------------------------------------------------------------------
$hdr1 = Find-ParaByText "Week beginning April 5th"
$r1p1 = $hdr1.Next()
$r1p2 = $r1p1.Next()
$r1p3 = $r1p2.Next()

$r1p3.Range.Delete()
$r1p2.Range.Text = "Initial team planning"
$r1p1.Range.Text = "Create Slack channel"

# ---------------------------------------------------------------------
# Row "Week beginning Monday April 12th": first trailing paragraph
# carries bold paragraph-mark formatting (<w:pPr><w:rPr><w:b/>...),
# the other two are plain empty paragraphs. End state: 2 plain
# paragraphs with notes text and no leftover bold paragraph formatting.
# ---------------------------------------------------------------------
$hdr2 = Find-ParaByText "April 12th"
$r2p1 = $hdr2.Next()
$r2p2 = $r2p1.Next()
$r2p3 = $r2p2.Next()

# Delete the 3rd (extra) empty paragraph first.
$r2p3.Range.Delete()

# Fill in the 2nd paragraph's text (plain paragraph, safe to set directly).
$r2p2.Range.Text = "Work on adding photo upload feature for notes"

# The 1st paragraph only carries bold formatting on its paragraph mark
# (no run to target with Font properties), so insert a fresh,
# unformatted paragraph ahead of it, populate that, and delete the old
# bold-formatted one.
$r2p1.Range.InsertParagraphBefore()
$hdr2 = Find-ParaByText "April 12th"
$newPara = $hdr2.Next()
$oldBoldPara = $newPara.Next()
$newPara.Range.Text = "Fill out project planning document with help of team"
$oldBoldPara.Range.Delete()
